$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 5, pushing existing rows 5-11 down to 7-13
$ws.Rows("5:6").Insert()

# New row 5: weekly price entry (Primera)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C5").Value = "Los Lagos"
$ws.Range("D5").Value = 44519
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = "Otros"
$ws.Range("I5").Value = 100107002
$ws.Range("J5").Value = "Chirimoya"
$ws.Range("K5").Value = "Cultivar IV Región"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 400
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("Q5").Value = "$/bandeja 8 kilos"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 2688
$ws.Range("T5").Value = 8

# New row 6: weekly price entry (Segunda)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C6").Value = "Los Lagos"
$ws.Range("D6").Value = 44519
$ws.Range("D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107002
$ws.Range("J6").Value = "Chirimoya"
$ws.Range("K6").Value = "Cultivar IV Región"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("Q6").Value = "$/bandeja 8 kilos"
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 2250
$ws.Range("T6").Value = 8
